$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1
$ws.Range("A1").Value = "Nén mp3"
$ws.Range("B1").Value = "EKOS MP3Minimizer"

# Row 2 - hyperlink to yoyogames resources
$ws.Range("A2").Value = "Các tập tin âm thanh (mp3, wav)"
$ws.Hyperlinks.Add($ws.Range("B2"), "http://www.yoyogames.com/resources?cat_id=4")

# Row 3
$ws.Range("B3").Value = "The Witcher 2 Bonus Disc"

# Row 4 - hyperlink to gamedev.net sprites thread
$ws.Range("A4").Value = "Các tập tin hình ảnh (jpg)"
$ws.Hyperlinks.Add($ws.Range("B4"), "http://www.gamedev.net/topic/272386-sprites-sprites-and-more-sprites/")

# Row 5
$ws.Range("B5").Value = "Bejeweled 2"

# Row 6
$ws.Range("B6").Value = "Resource cung cấp trong môn học C4W"

# Size the columns to fit their content, like Excel's column auto-fit.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
